# Adds basic webservice / blazer app authentication sample data to the
# TestData workbook: a "Domain" based external-user identifier column on
# the Users sheet, plus matching sample rows on Wallets and
# WalletTransactions, and makes "Users" the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Users sheet -----------------------------------------------------
$users = $wb.Worksheets.Item("Users")

# New sample user row (row 8) - entered first
$users.Cells.Item(8, 2).Value = "Domain"

# New column C header
$users.Cells.Item(1, 3).Value = "UniqueExternalUserID"

# GUID-looking external id stored as text (leading apostrophe forces the
# quote-prefix so Excel keeps it as text instead of re-interpreting it)
$users.Cells.Item(6, 3).Value = "'5f5a978a-4595-4e07-b456-73bf26fe6786"

$users.Cells.Item(8, 1).Value = 7

$users.Columns.Item(3).ColumnWidth = 36.140625

# --- Wallets sheet -----------------------------------------------------
$wallets = $wb.Worksheets.Item("Wallets")
$wallets.Cells.Item(8, 1).Value = 7
$wallets.Cells.Item(8, 2).Value = 7
$wallets.Cells.Item(8, 3).Value = 1000
$wallets.Cells.Item(8, 3).NumberFormat = $wallets.Cells.Item(7, 3).NumberFormat

# --- WalletTransactions sheet -------------------------------------------
$wtx = $wb.Worksheets.Item("WalletTransactions")
$wtx.Cells.Item(9, 1).Value = 8
$wtx.Cells.Item(9, 2).Value = 7
$wtx.Cells.Item(9, 3).Value = 1
$wtx.Cells.Item(9, 5).Value = 1000
$wtx.Cells.Item(9, 5).NumberFormat = $wtx.Cells.Item(8, 5).NumberFormat

# --- Make Users the active sheet/tab ------------------------------------
$users.Activate()
